$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.NumberFormat = "General"
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '45.764.05'
Set-TextValue $ws.Range('E2') '  +6.54%  '
Set-TextValue $ws.Range('D3') '2.392.80'
Set-TextValue $ws.Range('E3') '  +4.04%  '
Set-TextValue $ws.Range('E4') '  -0.10%  '
Set-TextValue $ws.Range('D5') '112.86'
Set-TextValue $ws.Range('E5') '  +7.55%  '
Set-TextValue $ws.Range('D6') '317.66'
Set-TextValue $ws.Range('E6') '  +2.40%  '
Set-TextValue $ws.Range('D7') '0.633'
Set-TextValue $ws.Range('E7') '  +1.14%  '
Set-TextValue $ws.Range('E8') '  -0.30%  '
Set-TextValue $ws.Range('D9') '0.625'
Set-TextValue $ws.Range('E9') '  +3.44%  '
Set-TextValue $ws.Range('D10') '41.68'
Set-TextValue $ws.Range('E10') '  +4.99%  '
Set-TextValue $ws.Range('D11') '0.0928'
Set-TextValue $ws.Range('E11') '  +2.48%  '
Set-TextValue $ws.Range('D12') '8.68'
Set-TextValue $ws.Range('E12') '  +4.76%  '
Set-TextValue $ws.Range('E13') '  +2.38%  '
Set-TextValue $ws.Range('D14') '1.00'
Set-TextValue $ws.Range('E14') '  +1.21%  '
Set-TextValue $ws.Range('D15') '15.76'
Set-TextValue $ws.Range('E15') '  +3.15%  '
Set-TextValue $ws.Range('D16') '2.760.09'
Set-TextValue $ws.Range('E16') '  +4.15%  '
Set-TextValue $ws.Range('D17') '2.407.65'
Set-TextValue $ws.Range('E17') '  +4.71%  '
Set-TextValue $ws.Range('D18') '45.702.62'
Set-TextValue $ws.Range('E18') '  +6.69%  '
Set-TextValue $ws.Range('D19') '7.46'
Set-TextValue $ws.Range('E19') '  +1.94%  '
Set-TextValue $ws.Range('D20') '0.0000108'
Set-TextValue $ws.Range('E20') '  +3.15%  '
Set-TextValue $ws.Range('D21') '13.34'
Set-TextValue $ws.Range('E21') '  -2.96%  '
Set-TextValue $ws.Range('D22') '74.40'
Set-TextValue $ws.Range('E22') '  +1.25%  '
Set-TextValue $ws.Range('D23') '3.52'
Set-TextValue $ws.Range('E23') '  +1.42%  '
Set-TextValue $ws.Range('D24') '264.17'
Set-TextValue $ws.Range('E24') '  -1.20%  '
Set-TextValue $ws.Range('D25') '2.33'
Set-TextValue $ws.Range('E25') '  +4.01%  '
Set-TextValue $ws.Range('E26') '  -0.70%  '
Set-TextValue $ws.Range('D27') '7.66'
Set-TextValue $ws.Range('E27') '  +1.11%  '
Set-TextValue $ws.Range('D28') '11.27'
Set-TextValue $ws.Range('E28') '  +2.98%  '
Set-TextValue $ws.Range('E29') '  +2.41%  '
Set-TextValue $ws.Range('D30') '38.84'
Set-TextValue $ws.Range('D31') '22.70'
Set-TextValue $ws.Range('E31') '  +2.22%  '
Set-TextValue $ws.Range('D32') '0.0976'
Set-TextValue $ws.Range('E32') '  +12.77%  '
Set-TextValue $ws.Range('D33') '171.98'
Set-TextValue $ws.Range('E33') '  +4.14%  '
Set-TextValue $ws.Range('E34') '  +4.18%  '
Set-TextValue $ws.Range('D35') '0.132'
Set-TextValue $ws.Range('E35') '  +1.30%  '
Set-TextValue $ws.Range('D36') '0.118'
Set-TextValue $ws.Range('E36') '  +5.29%  '
Set-TextValue $ws.Range('D37') '4.87'
Set-TextValue $ws.Range('E37') '  +5.33%  '
Set-TextValue $ws.Range('D38') '4.08'
Set-TextValue $ws.Range('E38') '  +12.92%  '
Set-TextValue $ws.Range('D39') '3.02'
Set-TextValue $ws.Range('E39') '  +7.87%  '
Set-TextValue $ws.Range('D40') '0.0361'
Set-TextValue $ws.Range('E40') '  +1.38%  '
Set-TextValue $ws.Range('D41') '1.77'
Set-TextValue $ws.Range('E41') '  +12.98%  '
Set-TextValue $ws.Range('D42') '102.39'
Set-TextValue $ws.Range('E42') '  -4.81%  '
Set-TextValue $ws.Range('D43') '0.239'
Set-TextValue $ws.Range('E43') '  +4.87%  '
Set-TextValue $ws.Range('D44') '13.42'
Set-TextValue $ws.Range('E44') '  +9.47%  '
Set-TextValue $ws.Range('D45') '71.45'
Set-TextValue $ws.Range('E45') '  +0.02%  '
Set-TextValue $ws.Range('D46') '87.50'
Set-TextValue $ws.Range('E46') '  +15.73%  '
Set-TextValue $ws.Range('E47') '  -0.58%  '
Set-TextValue $ws.Range('D48') '114.89'
Set-TextValue $ws.Range('E48') '  +2.90%  '
Set-TextValue $ws.Range('D49') '9.48'
Set-TextValue $ws.Range('E49') '  +6.95%  '
Set-TextValue $ws.Range('D50') '5.60'
Set-TextValue $ws.Range('E50') '  +7.99%  '
Set-TextValue $ws.Range('D51') '1.656.96'
Set-TextValue $ws.Range('E51') '  -2.93%  '
